# Mise à jour des résultats du script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" placeholder value in C162 -> becomes a blank cell
$ws.Range("C162").Value = $null

# New data rows appended to the table (rows 163-170)
$newRows = @(
    @("2025-07-04", "eaux souterraines", 115, 1),
    @("2025-07-04", "ruissellement",     116, 1),
    @("2025-07-04", "eaux souterraines", 116, 1),
    @("2025-07-04", "eaux souterraines", 117, 1),
    @("2025-07-04", "eaux de surface",   117, 1),
    @("2025-07-04", "ruissellement",     117, 1),
    @("2025-07-04", "eaux souterraines", 119, 1),
    @("2025-07-04", "eaux souterraines", 122, 1)
)

$startRow = 163
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date written as plain text (e.g. "2025-07-04"), matching
    # the rest of the sheet. Force text formatting first so Excel doesn't
    # auto-convert it to a date serial, then strip the format again so the
    # cell ends up styleless just like its neighbours.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
